# Rollerball System Test Cases.xlsx
#
# The only meaningful content change in the diff is to test case #8
# ("Make a move logout and login"), whose Test Steps cell (C9) had its
# final step text tightened up:
#   "7: go back to 1 and go through 3"
# becomes
#   "7: go back to step.1 and go through step.3"
#
# (All the other hunks in the diff are just the shared-strings table being
# reshuffled/renumbered as a side effect of that edit - the actual cell
# values elsewhere on the sheet are unchanged.)
#
# The sheet's active selection also moved from E14 to C9 (i.e. the author
# ended up with C9 selected after making the edit).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$newSteps = "1: Connect to the system" + "`n" + `
            "2: Log in " + "`n" + `
            "3: Select Games button" + "`n" + `
            "4: select  piece" + "`n" + `
            "5: Select a highlighted square" + "`n" + `
            "6: close Client" + "`n" + `
            "7: go back to step.1 and go through step.3"

$ws.Range("C9").Value = $newSteps

# Leave the workbook with C9 selected, matching the saved sheet view.
$ws.Range("C9").Select()
